$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Before touching anything, clone the still-blank placeholder row (595)
#    down into the new trailing placeholder rows (603-618) so they pick up
#    the exact same "empty" styling (date column style 4, text columns
#    style 5) that the rest of the blank padding rows already use.
# ---------------------------------------------------------------------------
$ws.Range("A595:E595").Copy()
$ws.Range("A603:E618").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Fill in the eight new event rows (595-602).
# ---------------------------------------------------------------------------
$ws.Range("A595").Value = 45850
$ws.Range("B595").Value = "AFTER PARTY"
$ws.Range("C595").Value = "SNRS"
$ws.Range("D595").Value = "Dortmund"
$ws.Range("E595").Value = "https://www.instagram.com/reel/DLkjOn-st9K/?igsh=MWoyMGp1NXE0NzBnYQ=="

$ws.Range("A596").Value = 45843
$ws.Range("B596").Value = "CLUB NIGHT"
$ws.Range("C596").Value = "SNRS"
$ws.Range("D596").Value = "Dortmund"
$ws.Range("E596").Value = "https://www.instagram.com/reel/DLZeJ2fIayI/?igsh=NHNqdmRvMGM3cTZn"

$ws.Range("A597").Value = 46003
$ws.Range("B597").Value = "UEBERREST"
$ws.Range("C597").Value = "SNRS"
$ws.Range("D597").Value = "Dortmund"
$ws.Range("E597").Value = "https://www.instagram.com/snrs.do?igsh=MXdzZTZtMzgyZG1zMg=="

$ws.Range("A598").Value = 45840
$ws.Range("B598").Value = "150 MPH RAVE"
$ws.Range("C598").Value = "Kanal"
$ws.Range("D598").Value = "Recklinghausen"
$ws.Range("E598").Value = "https://www.instagram.com/reel/DLjDHTtIU9z/?igsh=dG4yeWNwam90bmRy"

$ws.Range("A599").Value = 45843
$ws.Range("B599").Value = "AFTERPARTY"
$ws.Range("C599").Value = "Mikroport"
$ws.Range("D599").Value = "Krefeld"
$ws.Range("E599").Value = "https://www.instagram.com/p/DLjqp-BIHvn/?igsh=YnZ6MHZwb2Nma2Fz"

$ws.Range("A600").Value = 45843
$ws.Range("B600").Value = "RAVE GATE x TRANCEFLIGHT"
$ws.Range("C600").Value = "Sam" + [char]0x2018 + "s"
$ws.Range("D600").Value = "Bielefeld"
$ws.Range("E600").Value = "https://www.instagram.com/reel/DLatxckscuX/?igsh=MWZ3aTViaHZ1NzhnOA=="

$ws.Range("A601").Value = 45843
$ws.Range("B601").Value = "RAVE TOTAL"
$ws.Range("C601").Value = "AREA 15"
$ws.Range("D601").Value = "Bochum"
$ws.Range("E601").Value = "https://www.instagram.com/reel/DLiLAGLMMZY/?igsh=MTF3cmhpZmQ0b2txNg=="

$ws.Range("A602").Value = 45871
$ws.Range("B602").Value = "B2B SPECIAL ALAADDIN ALL NIGHT LONG"
$ws.Range("C602").Value = "Sam" + [char]0x2018 + "s"
$ws.Range("D602").Value = "Bielefeld"
$ws.Range("E602").Value = "https://www.instagram.com/reel/DLkjGQWM_lH/?igsh=b2M3cmhyejB0ejk4"

# ---------------------------------------------------------------------------
# 3) Re-apply the standard "filled row" formatting (style used by every
#    other populated data row) onto the new rows, since a plain .Value
#    assignment keeps whatever formatting the cell already had (the blank
#    placeholder look).
# ---------------------------------------------------------------------------
$ws.Range("A590:E590").Copy()
$ws.Range("A595:E602").PasteSpecial(-4122)
$ws.Range("A595:E618").RowHeight = 15

# ---------------------------------------------------------------------------
# 4) Wire up the hyperlinks on column E for the new rows, then restore the
#    plain "filled row" cell style (Hyperlinks.Add swaps in its own
#    hyperlink style) so the cell formatting matches the rest of the sheet.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("E595"), "https://www.instagram.com/reel/DLkjOn-st9K/?igsh=MWoyMGp1NXE0NzBnYQ==", "", "", "https://www.instagram.com/reel/DLkjOn-st9K/?igsh=MWoyMGp1NXE0NzBnYQ==")
$ws.Hyperlinks.Add($ws.Range("E596"), "https://www.instagram.com/reel/DLZeJ2fIayI/?igsh=NHNqdmRvMGM3cTZn", "", "", "https://www.instagram.com/reel/DLZeJ2fIayI/?igsh=NHNqdmRvMGM3cTZn")
$ws.Hyperlinks.Add($ws.Range("E597"), "https://www.instagram.com/snrs.do?igsh=MXdzZTZtMzgyZG1zMg==", "", "", "https://www.instagram.com/snrs.do?igsh=MXdzZTZtMzgyZG1zMg==")
$ws.Hyperlinks.Add($ws.Range("E598"), "https://www.instagram.com/reel/DLjDHTtIU9z/?igsh=dG4yeWNwam90bmRy", "", "", "https://www.instagram.com/reel/DLjDHTtIU9z/?igsh=dG4yeWNwam90bmRy")
$ws.Hyperlinks.Add($ws.Range("E599"), "https://www.instagram.com/p/DLjqp-BIHvn/?igsh=YnZ6MHZwb2Nma2Fz", "", "", "https://www.instagram.com/p/DLjqp-BIHvn/?igsh=YnZ6MHZwb2Nma2Fz")
$ws.Hyperlinks.Add($ws.Range("E600"), "https://www.instagram.com/reel/DLatxckscuX/?igsh=MWZ3aTViaHZ1NzhnOA==", "", "", "https://www.instagram.com/reel/DLatxckscuX/?igsh=MWZ3aTViaHZ1NzhnOA==")
$ws.Hyperlinks.Add($ws.Range("E601"), "https://www.instagram.com/reel/DLiLAGLMMZY/?igsh=MTF3cmhpZmQ0b2txNg==", "", "", "https://www.instagram.com/reel/DLiLAGLMMZY/?igsh=MTF3cmhpZmQ0b2txNg==")
$ws.Hyperlinks.Add($ws.Range("E602"), "https://www.instagram.com/reel/DLkjGQWM_lH/?igsh=b2M3cmhyejB0ejk4", "", "", "https://www.instagram.com/reel/DLkjGQWM_lH/?igsh=b2M3cmhyejB0ejk4")

$ws.Range("E590").Copy()
$ws.Range("E595:E602").PasteSpecial(-4122)
